$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Chart")

# Append the new day's data as row 84 (one row below the current last row, 83).
# Column A holds the date as plain text (matching every other row in the
# column), so we can't just assign the string via .Value -- Excel's normal
# "looks like a date" auto-detection would convert it into a real date
# serial number + date-formatted style. Instead we build the text through a
# formula that evaluates to the literal string, then paste-special just the
# value over the target cell; that carries the text across without
# re-running the "is this a date" input parser.
$scratch = $ws.Cells.Item(90, 1)
$scratch.Formula = "=""2025-12-28"""
$scratch.Copy()
$ws.Cells.Item(84, 1).PasteSpecial(-4163)
$scratch.ClearContents()

$ws.Cells.Item(84, 2).Value = 0.0
$ws.Cells.Item(84, 3).Value = 28.0
